$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.343.27"
$ws.Range("E2").Value = "  +2.23%  "
$ws.Range("D3").Value = "2.467.39"
$ws.Range("E3").Value = "  +2.36%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "574.62"
$ws.Range("E5").Value = "  +2.21%  "
$ws.Range("D6").Value = "148.74"
$ws.Range("E6").Value = "  +4.59%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("D8").Value = "0.538"
$ws.Range("E8").Value = "  +1.56%  "
$ws.Range("D9").Value = "0.113"
$ws.Range("E9").Value = "  +3.33%  "
$ws.Range("E10").Value = "  +0.68%  "
$ws.Range("D11").Value = "5.35"
$ws.Range("E11").Value = "  +2.53%  "
$ws.Range("D12").Value = "'0.360"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.14%  "
$ws.Range("D13").Value = "27.14"
$ws.Range("E13").Value = "  +6.20%  "
$ws.Range("D14").Value = "0.0000185"
$ws.Range("E14").Value = "  +6.70%  "
$ws.Range("D15").Value = "2.852.05"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").Value = "62.978.35"
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("D17").Value = "2.451.35"
$ws.Range("E17").Value = "  +1.39%  "
$ws.Range("D18").Value = "11.38"
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("D19").Value = "7.16"
$ws.Range("E19").Value = "  +5.21%  "
$ws.Range("D20").Value = "327.39"
$ws.Range("E20").Value = "  +1.78%  "
$ws.Range("D21").Value = "4.21"
$ws.Range("E21").Value = "  +1.89%  "
$ws.Range("D22").Value = "1.01"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").Value = "'1.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.50%  "
$ws.Range("D24").Value = "67.61"
$ws.Range("E24").Value = "  +2.01%  "
$ws.Range("D25").Value = "654.68"
$ws.Range("E25").Value = "  +16.67%  "
$ws.Range("D26").Value = "8.94"
$ws.Range("E26").Value = "  +2.75%  "
$ws.Range("D27").Value = "0.0000105"
$ws.Range("E27").Value = "  +12.41%  "
$ws.Range("D28").Value = "2.585.10"
$ws.Range("E28").Value = "  +2.07%  "
$ws.Range("D29").Value = "8.55"
$ws.Range("E29").Value = "  +4.64%  "
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "'1.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.43%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "0.996"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").Value = "0.144"
$ws.Range("E32").Value = "  -2.28%  "
$ws.Range("D33").Value = "'1.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("D34").Value = "1.54"
$ws.Range("E34").Value = "  +2.02%  "
$ws.Range("D35").Value = "'5.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.23%  "
$ws.Range("D36").Value = "0.997"
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("D37").Value = "0.386"
$ws.Range("E37").Value = "  +1.91%  "
$ws.Range("D38").Value = "5.53"
$ws.Range("E38").Value = "  +1.94%  "
$ws.Range("D39").Value = "18.87"
$ws.Range("E39").Value = "  +1.90%  "
$ws.Range("D40").Value = "1.87"
$ws.Range("E40").Value = "  +3.49%  "
$ws.Range("D41").Value = "148.49"
$ws.Range("E41").Value = "  -2.84%  "
$ws.Range("D42").Value = "2.63"
$ws.Range("E42").Value = "  +17.61%  "
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "'152.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.02%  "
$ws.Range("D45").Value = "3.73"
$ws.Range("E45").Value = "  +3.24%  "
$ws.Range("D46").Value = "0.0549"
$ws.Range("E46").Value = "  +3.83%  "
$ws.Range("D47").Value = "21.03"
$ws.Range("E47").Value = "  +6.15%  "
$ws.Range("D48").Value = "0.613"
$ws.Range("E48").Value = "  +3.44%  "
$ws.Range("D49").Value = "0.0236"
$ws.Range("E49").Value = "  +4.24%  "
$ws.Range("D50").Value = "'0.0930"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("E51").Value = "  +4.10%  "
